$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samu Lyhty")

# Update hours for row 19 (2026-02-01 entry) from 4 to 6
$ws.Range("C19").Value = 6

# Fill in new rows 28-32 with dates, tasks, hours
$ws.Range("A28").Value = 46069
$ws.Range("B28").Value = "Palaveri/Projektin viimeistelyä/videoiden suunnittelua"
$ws.Range("C28").Value = 3

$ws.Range("A29").Value = 46070
$ws.Range("B29").Value = "Projektin viimeistelyä"
$ws.Range("C29").Value = 5

$ws.Range("A30").Value = 46071
$ws.Range("B30").Value = "Esityksen tekoa/Projektin viimeistelyä"
$ws.Range("C30").Value = 4

$ws.Range("A31").Value = 46072
$ws.Range("B31").Value = "Esityksen tekoa/Projektin viimeistelyä"
$ws.Range("C31").Value = 7

$ws.Range("A32").Value = 46073
$ws.Range("B32").Value = "Sovelluksen julkaisemista/Projektin viimeistelyä"
$ws.Range("C32").Value = 5

# Row 33: total row (copy date-cell formatting from A32 so style index is reused, then set text)
$ws.Cells.Item(32, 1).Copy()
$ws.Cells.Item(33, 1).PasteSpecial(-4122)
$ws.Cells.Item(33, 1).Value = "Yht"
$ws.Range("C33").Formula = "=SUM(C3:C32)"

# Remove old SUM row at 100
$ws.Rows.Item(100).Delete()

# Update sheet view: selecting E16 resets scrolling so topLeftCell reverts to default
$ws.Range("E16").Select()
